$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [char]0x2083

$ws.Range("D2").Value = "58.526.35"
$ws.Range("E2").Value = "  -2.82%  "
$ws.Range("D3").Value = "2.303.44"
$ws.Range("E3").Value = "  -4.66%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.97%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.573"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.50%  "
$ws.Range("D9").Value = "2.301.28"
$ws.Range("E9").Value = "  -4.66%  "
$ws.Range("E10").Value = "  -4.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.150"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("E13").Value = "  -5.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.49%  "
$ws.Range("D15").Value = "2.713.28"
$ws.Range("E15").Value = "  -4.72%  "
$ws.Range("D16").Value = "58.517.21"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("E17").Value = "  -4.10%  "
$ws.Range("D18").Value = "2.313.37"
$ws.Range("E18").Value = "  -3.75%  "
$ws.Range("E19").Value = "  -5.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "312.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.13%  "
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.88%  "
$ws.Range("E25").Value = "  -4.00%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.82%  "
$ws.Range("E29").Value = "  -3.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").Value = "0.0${sub3}0719"
$ws.Range("E31").Value = "  -6.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.378"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  -3.79%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.95%  "
$ws.Range("E39").Value = "  -6.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "291.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "139.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0949"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.29%  "
$ws.Range("E46").Value = "  -3.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.563"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.30%  "
$ws.Range("E49").Value = "  -3.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.04%  "
